# Uppercase all textual data in the worksheet, except the header row (row 1).
# This mirrors the "GaloScrapper" commit: all scraped data values are converted
# to upper case for consistent display/manipulation, while column headers in
# row 1 (Title, Type, Author 1, ...) are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow  = $usedRange.Row
$firstCol  = $usedRange.Column
$lastRow   = $firstRow + $usedRange.Rows.Count - 1
$lastCol   = $firstCol + $usedRange.Columns.Count - 1

# Data starts on row 2; row 1 holds the column headers and must stay untouched.
$dataStartRow = [Math]::Max($firstRow + 1, 2)

for ($r = $dataStartRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -is [string]) {
            $upper = $val.ToUpper()
            if (-not $upper.Equals($val)) {
                $cell.Value = $upper
            }
        }
    }
}
